# Donor Application.xlsx — add an "ADDITIONAL INFORMATION" section
# (Org / Sports *) between CONTACT INFORMATION and SCHOLARSHIP TO OFFER
# DETAILS, per commit "feat: multiple uploads per one csv".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert 3 blank rows right before the old row 15
#        ("SCHOLARSHIP TO OFFER DETAILS"). Everything below shifts down
#        by 3 (dimension, merged cells and data-validation ranges all
#        follow automatically).
$ws.Rows("15:17").Insert()

# --- 2. New row 15: section header "ADDITIONAL INFORMATION".
#        Clone the look of the other section headers (row 12,
#        "CONTACT INFORMATION") so the style indices are reused as-is.
$ws.Range("A12:D12").Copy()
$ws.Range("A15:D15").PasteSpecial(-4122)
$ws.Range("B15").Value = "ADDITIONAL INFORMATION"
$ws.Rows(15).RowHeight = 14.5
$ws.Range("B15:C15").Merge()

# --- 3. New row 16: "Org" field label, styled like the other plain
#        field rows (row 13, "Contact Number *").
$ws.Range("A13:D13").Copy()
$ws.Range("A16:D16").PasteSpecial(-4122)
$ws.Range("B16").Value = "Org"
$ws.Rows(16).RowHeight = 14.5

# --- 4. New row 17: "Sports *" field label — last row of the new
#        section, slightly taller (matches the row just above every
#        other section header in this sheet).
$ws.Range("A13:D13").Copy()
$ws.Range("A17:D17").PasteSpecial(-4122)
$ws.Range("B17").Value = "Sports *"
$ws.Rows(17).RowHeight = 15

# --- 5. The rows that now sit right above a section header (row 11,
#        "Citizenship *", above the existing CONTACT INFORMATION header,
#        and row 14, "Email *", above the new ADDITIONAL INFORMATION
#        header) pick up the same slightly-taller row height.
$ws.Rows(11).RowHeight = 15
$ws.Rows(14).RowHeight = 15

# --- 6. Move the active selection the same amount the content below
#        it shifted.
$ws.Range("C17").Select()

Write-Host "Added ADDITIONAL INFORMATION section (Org / Sports *) at rows 15-17"
